# Convert the numeric 0/1 values in column H (rows 2-114) of the
# "pkm_conquest_pokemon_evolution" sheet into the text strings "FALSE" / "TRUE"
# (stored as shared-string cells, t="s"), matching the target OOXML diff.
#
# A plain `$ws.Range(...).Value = "FALSE"` gets auto-coerced to an Excel
# Boolean (t="b") by the object model, and prefixing with a leading quote
# forces text but adds a spurious quotePrefix cell style. Routing the text
# through a formula ("="..."") and then collapsing it to a plain value via
# Copy + PasteSpecial(xlPasteValues) yields a clean shared-string cell with
# no extra style - exactly what the diff expects.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

# Rows (in column H) that must end up as the text "TRUE"; all remaining
# data rows (2-114) become "FALSE".
$trueRows = @(12, 14, 16, 32, 33, 68, 72, 86, 89)

# First, stamp the whole H2:H114 block with "FALSE" in one shot (fast path -
# Copy/PasteSpecial behaves correctly over a single contiguous range).
$falseRange = $ws.Range("H2:H114")
$falseRange.Formula = "=""FALSE"""
$falseRange.Copy()
$falseRange.PasteSpecial($xlPasteValues)

# Then overwrite just the TRUE rows individually (a multi-area Union only
# round-trips its first area through Copy/PasteSpecial reliably here, so we
# touch each TRUE cell on its own).
foreach ($r in $trueRows) {
    $cell = $ws.Range("H$r")
    $cell.Formula = "=""TRUE"""
    $cell.Copy()
    $cell.PasteSpecial($xlPasteValues)
}
